$p = $ppt.ActivePresentation

# Slide 3 ("An Illustrative Example" / chi-square rows discussion) contains the
# bullet text that originally had the sentence split across two runs:
#   "...then the total row could represent either row" + "."
# It is merged back into a single run/sentence:
#   "...then the total row could represent either row."
$targetSlideIndex = 3
$targetShapeName = "Rectangle 5"
$needle = "the total row could represent either row."

$slide = $p.Slides.Item($targetSlideIndex)

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)

    if ($shape.HasTextFrame -and $shape.Name -eq $targetShapeName) {
        $textRange = $shape.TextFrame.TextRange
        $fullText = $textRange.Text

        if ($fullText.Contains($needle)) {
            $startPos = $fullText.IndexOf($needle) + 1
            $len = $needle.Length

            # Re-assigning the characters spanning the old run boundary causes
            # the (identically-formatted) adjacent runs to coalesce into one
            # run, exactly as PowerPoint does when the text is retyped/edited.
            $chars = $textRange.Characters($startPos, $len)
            $chars.Text = $needle
        }
    }
}
